# updated GSC export data
# Append three more days (2025-11-11, 2025-11-12, 2025-11-13) to the
# "Chart" sheet's daily video-indexing table, continuing the existing
# pattern (No video indexed = 24, Video indexed = 0, Impressions = 0).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chart")

$newDates = @("2025-11-11", "2025-11-12", "2025-11-13")

$lastRow = 38
$row = $lastRow + 1
foreach ($d in $newDates) {
    $dateCell = $ws.Cells.Item($row, 1)
    # Force the date-looking string to be stored as text (matching the
    # existing rows, which keep dates as shared strings, not date
    # serials), then drop the temporary text format so the cell keeps
    # the sheet's default style.
    $dateCell.NumberFormat = "@"
    $dateCell.Value = $d
    $dateCell.ClearFormats()

    $ws.Cells.Item($row, 2).Value = 24
    $ws.Cells.Item($row, 3).Value = 0
    $ws.Cells.Item($row, 4).Value = 0

    $row = $row + 1
}
